# Add Policy Directive #50 as a new row (row 51) on the "directives" sheet,
# matching the existing table's layout: A=id, B=date, C=source, D=keywords,
# F=title, G=content (column E left blank, as in other rows referencing
# "see adoc file").

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A51").Value = 50

# Set the number format before assigning the value so the date is stored
# using the custom mm-dd-yy format (matching PD #50's row) rather than
# triggering an auto-detected date format.
$ws.Range("B51").NumberFormat = "mm-dd-yy"
$ws.Range("B51").Value = "03/26/2021"

$ws.Range("C51").Value = "TC"
$ws.Range("D51").Value = "standard, template, collaboration"
$ws.Range("F51").Value = "OGC API coordination"
$ws.Range("G51").Value = "see adoc file"

# Match the wrap-text formatting used throughout the table, applied per
# populated cell so we don't create a stray empty cell in column E.
$ws.Range("A51").WrapText = $true
$ws.Range("B51").WrapText = $true
$ws.Range("C51").WrapText = $true
$ws.Range("D51").WrapText = $true
$ws.Range("F51").WrapText = $true
$ws.Range("G51").WrapText = $true

$ws.Rows.Item(51).RowHeight = 34

# Move the view/selection the way Excel would after entering data on row 51
# and pressing Enter in column G (lands on G52 of the next row).
$ws.Range("G52").Select()
